$wb = $excel.ActiveWorkbook

# --- Tabelle1 (sheet 1): add Dezimal/Hex columns and a new "UhrModul" row ---
$ws1 = $wb.Worksheets.Item("Tabelle1")

$ws1.Range("D28").Value = "Dezimal"
$ws1.Range("E28").Value = "Hex"

$ws1.Range("D29").Value = 32
$ws1.Range("E29").Value = "0x20"

$ws1.Range("D30").Value = 33
$ws1.Range("E30").Value = "0x21"

$ws1.Range("D31").Value = 34
$ws1.Range("E31").Value = "0x22"

$ws1.Range("D32").Value = 35
$ws1.Range("E32").Value = "0x23"

$ws1.Range("D33").Value = 36
$ws1.Range("E33").Value = "0x24"

$ws1.Range("D34").Value = 63
$ws1.Range("E34").Value = "0x3F"

# Insert a new row for the "UhrModul" (clock module) entry right after the
# "LCD Display" row, pushing everything below down by one.
[void]$ws1.Rows.Item(35).Insert()

$ws1.Range("A35").Value = "UhrModul"
$ws1.Range("C35").Value = "?"
$ws1.Range("E35").Value = "0x57 und 0x68 zumindest laut i2c Scanner"

# --- Active sheet / selection bookkeeping ---
# Tabelle3's selection moves to A29 (it is no longer the visible/active tab).
$ws3 = $wb.Worksheets.Item("Tabelle3")
[void]$ws3.Activate()
[void]$ws3.Range("A29").Select()

# Tabelle1 becomes the active tab (was Tabelle3), scrolled so row 13 is at
# the top, with I32 selected.
[void]$ws1.Activate()
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws1.Range("I32").Select()
